# Adds "(SourceN)" data-source suffixes to the header row of the stream
# join report and widens the affected columns so the longer headers stay
# readable (mirrors the original author widening C:E and F:H after the
# header text grew).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: note which data source each column came from ---------
$ws.Range("C2").Value = "FIRST_NAME (Source1)"
$ws.Range("D2").Value = "LAST_NAME (Source1)"
$ws.Range("E2").Value = "AGE (Source1)"
$ws.Range("F2").Value = "FIRST_NAME (Source2)"
$ws.Range("G2").Value = "LAST_NAME (Source2)"
$ws.Range("H2").Value = "AGE (Source2)"

# --- Column widths: widen to fit the new, longer header text -----------
# C/F = FIRST_NAME (Source#) -> widest
# D/G = LAST_NAME (Source#)  -> second widest
# E/H = AGE (Source#)        -> narrowest
$ws.Columns.Item(3).ColumnWidth = 25.333333333333332
$ws.Columns.Item(4).ColumnWidth = 24.833333333333332
$ws.Columns.Item(5).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = 25.333333333333332
$ws.Columns.Item(7).ColumnWidth = 24.833333333333332
$ws.Columns.Item(8).ColumnWidth = 17.5
